$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "67.710.33"
Set-TextValue $ws.Range("E2") "  -1.16%  "
Set-TextValue $ws.Range("D3") "3.797.76"
Set-TextValue $ws.Range("E3") "  +1.39%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.03%  "
Set-TextValue $ws.Range("D5") "595.46"
Set-TextValue $ws.Range("D6") "166.69"
Set-TextValue $ws.Range("E6") "  -0.55%  "
Set-TextValue $ws.Range("D7") "3.796.27"
Set-TextValue $ws.Range("E7") "  +1.40%  "
Set-TextValue $ws.Range("E9") "  -0.19%  "
Set-TextValue $ws.Range("D10") "0.160"
Set-TextValue $ws.Range("E10") "  -0.23%  "
Set-TextValue $ws.Range("E11") "  -1.86%  "
Set-TextValue $ws.Range("E12") "  +0.20%  "
Set-TextValue $ws.Range("D13") "0.0000257"
Set-TextValue $ws.Range("E13") "  -0.97%  "
Set-TextValue $ws.Range("D14") "36.30"
Set-TextValue $ws.Range("E14") "  -0.14%  "
Set-TextValue $ws.Range("D15") "4.430.52"
Set-TextValue $ws.Range("E15") "  +1.32%  "
Set-TextValue $ws.Range("D16") "3.771.67"
Set-TextValue $ws.Range("E16") "  +0.75%  "
Set-TextValue $ws.Range("D17") "18.67"
Set-TextValue $ws.Range("E17") "  +4.09%  "
Set-TextValue $ws.Range("D18") "67.687.38"
Set-TextValue $ws.Range("E18") "  -1.14%  "
Set-TextValue $ws.Range("E19") "  +0.34%  "
Set-TextValue $ws.Range("D20") "6.99"
Set-TextValue $ws.Range("E20") "  -0.42%  "
Set-TextValue $ws.Range("E21") "  -4.37%  "
Set-TextValue $ws.Range("D22") "459.06"
Set-TextValue $ws.Range("E22") "  -1.84%  "
Set-TextValue $ws.Range("D23") "0.699"
Set-TextValue $ws.Range("E23") "  -0.46%  "
Set-TextValue $ws.Range("D24") "0.0000154"
Set-TextValue $ws.Range("E24") "  +5.95%  "
Set-TextValue $ws.Range("D25") "83.75"
Set-TextValue $ws.Range("E25") "  -0.54%  "
Set-TextValue $ws.Range("D26") "11.93"
Set-TextValue $ws.Range("E26") "  -1.27%  "
Set-TextValue $ws.Range("D27") "2.14"
Set-TextValue $ws.Range("E27") "  -2.77%  "
Set-TextValue $ws.Range("E28") "  -0.58%  "
Set-TextValue $ws.Range("E30") "  +0.07%  "
Set-TextValue $ws.Range("E31") "  -0.54%  "
Set-TextValue $ws.Range("D32") "29.88"
Set-TextValue $ws.Range("E32") "  -0.29%  "
Set-TextValue $ws.Range("E33") "  +0.42%  "
Set-TextValue $ws.Range("E34") "  -0.78%  "
Set-TextValue $ws.Range("E35") "  -0.21%  "
Set-TextValue $ws.Range("D36") "3.747.81"
Set-TextValue $ws.Range("E36") "  +1.27%  "
Set-TextValue $ws.Range("E37") "  -1.50%  "
Set-TextValue $ws.Range("D38") "3.35"
Set-TextValue $ws.Range("E38") "  -3.48%  "
Set-TextValue $ws.Range("E39") "  -0.33%  "
Set-TextValue $ws.Range("D40") "1.00"
Set-TextValue $ws.Range("E40") "  +0.01%  "
Set-TextValue $ws.Range("D41") "5.76"
Set-TextValue $ws.Range("E41") "  -0.92%  "
Set-TextValue $ws.Range("E42") "  -0.02%  "
Set-TextValue $ws.Range("D44") "44.85"
Set-TextValue $ws.Range("E44") "  +4.24%  "
Set-TextValue $ws.Range("E45") "  -2.19%  "
Set-TextValue $ws.Range("E46") "  +2.65%  "
Set-TextValue $ws.Range("D47") "8.39"
Set-TextValue $ws.Range("E47") "  -2.72%  "
Set-TextValue $ws.Range("D48") "148.03"
Set-TextValue $ws.Range("E48") "  +0.79%  "
Set-TextValue $ws.Range("D49") "394.31"
Set-TextValue $ws.Range("E49") "  +0.33%  "
Set-TextValue $ws.Range("E50") "  -5.21%  "
Set-TextValue $ws.Range("D51") "2.759.88"
Set-TextValue $ws.Range("E51") "  +2.47%  "
